# Add a new "14-ago" date column (AY) to the tracking sheet, mirroring the
# existing AX ("11-ago") column, and populate it with the latest figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column.
$ws.Range("AY1").Value = "14-ago"

# New daily values, one per data row.
$ws.Range("AY2").Value = 0
$ws.Range("AY3").Value = 17.471036026470173
$ws.Range("AY4").Value = 18.62859577367583
$ws.Range("AY5").Value = 13.00852083249233
$ws.Range("AY6").Value = 0
$ws.Range("AY7").Value = 22.312450829281918
$ws.Range("AY8").Value = 12.909476434487644
$ws.Range("AY9").Value = 13.622743317224776
$ws.Range("AY10").Value = 10.379891936413133
$ws.Range("AY11").Value = 8.5388710638244518
$ws.Range("AY12").Value = 0
$ws.Range("AY13").Value = 13.695878812869779
$ws.Range("AY14").Value = 0
$ws.Range("AY15").Value = 0
$ws.Range("AY16").Value = 15.711033388543576
$ws.Range("AY17").Value = 0
$ws.Range("AY18").Value = 0

# Match the author's final selection in the sheet (the whole new column's data).
[void]$ws.Range("AY2:AY18").Select()
